$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 74197.169559412549
$ws.Range("B4").Value = 98202.13618157545
$ws.Range("B5").Value = 4364.5393858477964
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 50919.626168224291
$ws.Range("B8").Value = 3637.1161548731639

$ws.Range("B7:B8").Select()
